$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (ID) stays text-typed so leading zeros are preserved,
# matching the original workbook where the ID column is shared-string typed.
$ws.Range("D2:D31").NumberFormat = "@"

$ws.Range("A2").Value = "Examtaker42276"
$ws.Range("B2").Value = "Automation42276"
$ws.Range("C2").Value = "examtakerautomation42276@gmail.com"
$ws.Range("D2").Value = "42276"
$ws.Range("A3").Value = "Examtaker25213"
$ws.Range("B3").Value = "Automation25213"
$ws.Range("C3").Value = "examtakerautomation25213@gmail.com"
$ws.Range("D3").Value = "25213"
$ws.Range("A4").Value = "Examtaker34152"
$ws.Range("B4").Value = "Automation34152"
$ws.Range("C4").Value = "examtakerautomation34152@gmail.com"
$ws.Range("D4").Value = "34152"
$ws.Range("A5").Value = "Examtaker58043"
$ws.Range("B5").Value = "Automation58043"
$ws.Range("C5").Value = "examtakerautomation58043@gmail.com"
$ws.Range("D5").Value = "58043"
$ws.Range("A6").Value = "Examtaker16009"
$ws.Range("B6").Value = "Automation16009"
$ws.Range("C6").Value = "examtakerautomation16009@gmail.com"
$ws.Range("D6").Value = "16009"
$ws.Range("A7").Value = "Examtaker89714"
$ws.Range("B7").Value = "Automation89714"
$ws.Range("C7").Value = "examtakerautomation89714@gmail.com"
$ws.Range("D7").Value = "89714"
$ws.Range("A8").Value = "Examtaker53506"
$ws.Range("B8").Value = "Automation53506"
$ws.Range("C8").Value = "examtakerautomation53506@gmail.com"
$ws.Range("D8").Value = "53506"
$ws.Range("A9").Value = "Examtaker30412"
$ws.Range("B9").Value = "Automation30412"
$ws.Range("C9").Value = "examtakerautomation30412@gmail.com"
$ws.Range("D9").Value = "30412"
$ws.Range("A10").Value = "Examtaker21414"
$ws.Range("B10").Value = "Automation21414"
$ws.Range("C10").Value = "examtakerautomation21414@gmail.com"
$ws.Range("D10").Value = "21414"
$ws.Range("A11").Value = "Examtaker72256"
$ws.Range("B11").Value = "Automation72256"
$ws.Range("C11").Value = "examtakerautomation72256@gmail.com"
$ws.Range("D11").Value = "72256"
$ws.Range("A12").Value = "Examtaker78165"
$ws.Range("B12").Value = "Automation78165"
$ws.Range("C12").Value = "examtakerautomation78165@gmail.com"
$ws.Range("D12").Value = "78165"
$ws.Range("A13").Value = "Examtaker77090"
$ws.Range("B13").Value = "Automation77090"
$ws.Range("C13").Value = "examtakerautomation77090@gmail.com"
$ws.Range("D13").Value = "77090"
$ws.Range("A14").Value = "Examtaker41552"
$ws.Range("B14").Value = "Automation41552"
$ws.Range("C14").Value = "examtakerautomation41552@gmail.com"
$ws.Range("D14").Value = "41552"
$ws.Range("A15").Value = "Examtaker73759"
$ws.Range("B15").Value = "Automation73759"
$ws.Range("C15").Value = "examtakerautomation73759@gmail.com"
$ws.Range("D15").Value = "73759"
$ws.Range("A16").Value = "Examtaker81655"
$ws.Range("B16").Value = "Automation81655"
$ws.Range("C16").Value = "examtakerautomation81655@gmail.com"
$ws.Range("D16").Value = "81655"
$ws.Range("A17").Value = "Examtaker05997"
$ws.Range("B17").Value = "Automation05997"
$ws.Range("C17").Value = "examtakerautomation05997@gmail.com"
$ws.Range("D17").Value = "05997"
$ws.Range("A18").Value = "Examtaker52492"
$ws.Range("B18").Value = "Automation52492"
$ws.Range("C18").Value = "examtakerautomation52492@gmail.com"
$ws.Range("D18").Value = "52492"
$ws.Range("A19").Value = "Examtaker72373"
$ws.Range("B19").Value = "Automation72373"
$ws.Range("C19").Value = "examtakerautomation72373@gmail.com"
$ws.Range("D19").Value = "72373"
$ws.Range("A20").Value = "Examtaker68417"
$ws.Range("B20").Value = "Automation68417"
$ws.Range("C20").Value = "examtakerautomation68417@gmail.com"
$ws.Range("D20").Value = "68417"
$ws.Range("A21").Value = "Examtaker46605"
$ws.Range("B21").Value = "Automation46605"
$ws.Range("C21").Value = "examtakerautomation46605@gmail.com"
$ws.Range("D21").Value = "46605"
$ws.Range("A22").Value = "Examtaker25272"
$ws.Range("B22").Value = "Automation25272"
$ws.Range("C22").Value = "examtakerautomation25272@gmail.com"
$ws.Range("D22").Value = "25272"
$ws.Range("A23").Value = "Examtaker19611"
$ws.Range("B23").Value = "Automation19611"
$ws.Range("C23").Value = "examtakerautomation19611@gmail.com"
$ws.Range("D23").Value = "19611"
$ws.Range("A24").Value = "Examtaker59754"
$ws.Range("B24").Value = "Automation59754"
$ws.Range("C24").Value = "examtakerautomation59754@gmail.com"
$ws.Range("D24").Value = "59754"
$ws.Range("A25").Value = "Examtaker43968"
$ws.Range("B25").Value = "Automation43968"
$ws.Range("C25").Value = "examtakerautomation43968@gmail.com"
$ws.Range("D25").Value = "43968"
$ws.Range("A26").Value = "Examtaker07054"
$ws.Range("B26").Value = "Automation07054"
$ws.Range("C26").Value = "examtakerautomation07054@gmail.com"
$ws.Range("D26").Value = "07054"
$ws.Range("A27").Value = "Examtaker92122"
$ws.Range("B27").Value = "Automation92122"
$ws.Range("C27").Value = "examtakerautomation92122@gmail.com"
$ws.Range("D27").Value = "92122"
$ws.Range("A28").Value = "Examtaker80371"
$ws.Range("B28").Value = "Automation80371"
$ws.Range("C28").Value = "examtakerautomation80371@gmail.com"
$ws.Range("D28").Value = "80371"
$ws.Range("A29").Value = "Examtaker96876"
$ws.Range("B29").Value = "Automation96876"
$ws.Range("C29").Value = "examtakerautomation96876@gmail.com"
$ws.Range("D29").Value = "96876"
$ws.Range("A30").Value = "Examtaker13906"
$ws.Range("B30").Value = "Automation13906"
$ws.Range("C30").Value = "examtakerautomation13906@gmail.com"
$ws.Range("D30").Value = "13906"
$ws.Range("A31").Value = "Examtaker92035"
$ws.Range("B31").Value = "Automation92035"
$ws.Range("C31").Value = "examtakerautomation92035@gmail.com"
$ws.Range("D31").Value = "92035"

# Restore the original (default) cell style on column D after forcing text
# number-format, so the saved style index matches the rest of the data rows.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("D2:D31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
